$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = "Toncoin"
$ws.Range("B28").Value = "Cosmos"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D2").Value = "41.746.81"
$ws.Range("D3").Value = "2.285.71"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0903"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "2.629.56"
$ws.Range("D17").Value = "2.290.04"
$ws.Range("D18").Value = "41.703.84"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "279.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "163.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0868"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0992"
$ws.Range("D51").Style = "Normal"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("E6").Value = "  -6.05%  "
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("E10").Value = "  -7.91%  "
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("E12").Value = "  -4.33%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  -5.26%  "
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("E23").Value = "  +8.38%  "
$ws.Range("E24").Value = "  +6.40%  "
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  +6.59%  "
$ws.Range("E28").Value = "  -7.20%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("E31").Value = "  -6.26%  "
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("E34").Value = "  -4.78%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  -4.81%  "
$ws.Range("E37").Value = "  -3.94%  "
$ws.Range("E38").Value = "  +6.96%  "
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("E40").Value = "  -9.01%  "
$ws.Range("E41").Value = "  +15.78%  "
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  -8.43%  "
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("E47").Value = "  -4.49%  "
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("E50").Value = "  -6.78%  "
$ws.Range("E51").Value = "  -2.45%  "
